$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to host a "Fecha:" date stamp line,
# pushing the existing header/body rows down by one.
$ws.Rows(2).Insert()

# Label, right aligned, using the same bold "Aparajita" heading font
# used elsewhere in the report.
$ws.Range("F2").Value = "Fecha:"
$ws.Range("F2").Font.Name = "Aparajita"
$ws.Range("F2").Font.Bold = $true
$ws.Range("F2").HorizontalAlignment = -4152   # xlRight

# Today's date, plain "Aparajita" font, short-date format.
$ws.Range("G2").Font.Name = "Aparajita"
$ws.Range("G2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Formula = "=TODAY()"
$ws.Range("G2").HorizontalAlignment = -4131   # xlLeft

$ws.Range("G8").Select()
